$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 323, pushing the existing rows 323:341 down to 324:342.
$ws.Rows.Item(323).Insert()

# Populate the newly inserted row 323 with the new weekly price record.
$ws.Cells.Item(323, 1).Value = 5
$ws.Cells.Item(323, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(323, 3).Value = "Maule"
$ws.Cells.Item(323, 4).Value = 44706
$ws.Cells.Item(323, 5).Value = 7
$ws.Cells.Item(323, 6).Value = 100114013
$ws.Cells.Item(323, 7).Value = "Zanahoria"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 500
$ws.Cells.Item(323, 11).Value = 5500
$ws.Cells.Item(323, 12).Value = 5500
$ws.Cells.Item(323, 13).Value = 5500
$ws.Cells.Item(323, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(323, 15).Value = "Región de Ñuble"
$ws.Cells.Item(323, 16).Value = 275
$ws.Cells.Item(323, 17).Value = 20
$ws.Cells.Item(323, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the rest of column D.
$ws.Cells.Item(323, 4).NumberFormat = $ws.Cells.Item(324, 4).NumberFormat
